$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.386.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("E3").Value = "  +1.73%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").Value = "'168.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.08%  "

$ws.Range("D7").Value = "'3.805.70"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.70%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("E10").Value = "  +0.98%  "

$ws.Range("E11").Value = "  -0.93%  "

$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("D13").Value = "'0.0000261"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").Value = "'36.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.10%  "

$ws.Range("D15").Value = "'4.445.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.66%  "

$ws.Range("D16").Value = "'3.804.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("D17").Value = "'68.441.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.71%  "

$ws.Range("D18").Value = "'18.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.60%  "

$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("D21").Value = "'11.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.24%  "

$ws.Range("D22").Value = "'466.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("E24").Value = "  +8.92%  "

$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("D27").Value = "'11.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.26%  "

$ws.Range("E28").Value = "  -0.34%  "

$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("E30").Value = "  -0.68%  "

$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("D32").Value = "'30.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.27%  "

$ws.Range("E33").Value = "  -2.95%  "

$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("D36").Value = "'3.761.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.74%  "

$ws.Range("E37").Value = "  -0.62%  "

$ws.Range("D38").Value = "'3.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.21%  "

$ws.Range("E39").Value = "  +1.08%  "

$ws.Range("E40").Value = "  +1.22%  "

$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("D44").Value = "'44.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.07%  "

$ws.Range("D45").Value = "'0.301"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.08%  "

$ws.Range("D46").Value = "'47.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.17%  "

$ws.Range("D47").Value = "'1.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("D48").Value = "'8.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.75%  "

$ws.Range("D49").Value = "'396.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("D50").Value = "'146.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.29%  "

$ws.Range("D51").Value = "'2.809.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.77%  "
